$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("F3").Value = 525000000

# Row 4
$ws.Range("F4").Value = 25000000

# Row 5
$ws.Range("F5").Value = 50000000

# Row 6
$ws.Range("A6").Value = 95585
$ws.Range("B6").Value = "ZEGA"
$ws.Range("C6").Value = "MARGUERITE LOUIS"
$ws.Range("D6").Value = 58
$ws.Range("E6").Value = 102030405
$ws.Range("F6").Value = 70000000

# Row 7
$ws.Range("A7").Value = 64258
$ws.Range("B7").Value = "KOUDOU "
$ws.Range("C7").Value = "LAURENT"
$ws.Range("D7").Value = 73
$ws.Range("E7").Value = 1234569
$ws.Range("F7").Value = 65000
